$d = $word.ActiveDocument

# InsertXML only behaves as an in-place "replace" when the target Range spans
# the ENTIRE paragraph it belongs to (start == paragraph start, end == paragraph
# end); for any narrower range the new content gets appended at the end of the
# paragraph instead of at the range's position. So every edit below is expressed
# as a full-paragraph replacement, reproducing the paragraph's pPr and any
# unchanged trailing runs/bookmarks verbatim.
function Replace-ParagraphXml {
    param(
        [int]$paraIndex,
        [string]$newParaInnerXml
    )
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body><w:p>' + $newParaInnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1. "Khelil Bethel" -> "Khelil" (spellStart/spellEnd) + " Bethel"
# ---------------------------------------------------------------------------
$para2 = '<w:proofErr w:type="spellStart"/>' + '<w:r><w:t>Khelil</w:t></w:r>' + '<w:proofErr w:type="spellEnd"/>' + '<w:r><w:t xml:space="preserve"> Bethel</w:t></w:r>'
Replace-ParagraphXml 2 $para2

# ---------------------------------------------------------------------------
# 2. "NeoEpoch" heading -> wrapped in spellStart/spellEnd, same bold/size rPr
# ---------------------------------------------------------------------------
$titleRPr = '<w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$para4 = '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + '<w:proofErr w:type="spellStart"/>' + '<w:r>' + $titleRPr + '<w:t>NeoEpoch</w:t></w:r>' + '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml 4 $para4

# ---------------------------------------------------------------------------
# 3. "Our application..." paragraph -> split "collected all of the pickups"
# ---------------------------------------------------------------------------
$para5 = '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' + '<w:r w:rsidRPr="00623BAE"><w:t>Our application</w:t></w:r>' + '<w:r><w:t xml:space="preserve"> starts you off in a title screen which gives the user the option to either start the game, edit setting (notional for now), or quit back to the desktop. From there, if the user selects &#8220;start game&#8221; the game will begin. It begins in an overhead 2D world in which the player (using a placeholder sprite as neither of us are any good at art). The player must then explore the world and find people to start a civilization. The people are currently jewels as placeholders (resembling typical pickups of other games). Once the player has explored the entire world and collected </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>all of</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '<w:r><w:t xml:space="preserve"> the pickups, the player has completed the game and is sent back to the menu screen. The player must survive, as periodically he will be attacked by tigers, which are present in a 3D world combat scenario. We realized our original ideas were incredibly optimistic given the short time frame and presence of other classes/obligations. I do think we were able to represent the core concepts of our game, being survival and stone age.</w:t></w:r>'
Replace-ParagraphXml 5 $para5

# ---------------------------------------------------------------------------
# 4. "To the application..." paragraph -> several gramStart/gramEnd splits
#    plus "clicking" -> "left " + "clicking. "
# ---------------------------------------------------------------------------
$para7 = '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' + '<w:r><w:t xml:space="preserve">To the application, it is </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>pretty straight</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '<w:r><w:t xml:space="preserve"> forward for anyone who has played videos games before as the controls a pretty reminiscent of them. To select and item in the title-menu, the user simply clicks on the desired selection with their mouse. In the 2D overworld, the player can move around using w, a, s, d as well as the arrow keys. Picking up items is done automatically as the user walks over them. Once the player has encountered a battle, w, a, s, d </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>are</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '<w:r><w:t xml:space="preserve"> used for movement, the player can look around using the mouse, jumping is done by pressing the space bar, and attacking is done by </w:t></w:r>' + '<w:r><w:t xml:space="preserve">left </w:t></w:r>' + '<w:r><w:t xml:space="preserve">clicking. </w:t></w:r>'
Replace-ParagraphXml 7 $para7

# ---------------------------------------------------------------------------
# 5. "The final design..." paragraph -> gramStart/gramEnd + spellStart/spellEnd
#    splits, "challenged" -> "challenges", "Kehlil" -> "Khelil" (split "Khe"/"lil")
#    plus the trailing "." run and _GoBack bookmark are reproduced unchanged.
# ---------------------------------------------------------------------------
$para9 = '<w:pPr><w:ind w:firstLine="720"/></w:pPr>' + '<w:r><w:t xml:space="preserve">The final design of our game is an interesting compliment of 2D and 3D worlds, giving the player a wide variety of experiences. The map is designed in a way that it is aesthetically appealing, like an untainted world with no industrialization. Lots of greens, cliffs used as natural barriers, and water to give it a brighter feel. The map is set up in spokes so the player must visit each corner of the map to collect all the gems and complete the game. Combat is </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>pretty straight</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '<w:r><w:t xml:space="preserve"> forward, as previously mentioned, with an open field with just the player and the tiger for fighting. The player has a spear they can use, as it is in the &#8220;caveman&#8221; </w:t></w:r>' + '<w:r><w:t>ages. Major technical challenges</w:t></w:r>' + '<w:r><w:t xml:space="preserve"> were both of us using Unity for the first time and struggling to accomplish the </w:t></w:r>' + '<w:proofErr w:type="gramStart"/>' + '<w:r><w:t>most simple</w:t></w:r>' + '<w:proofErr w:type="gramEnd"/>' + '<w:r><w:t xml:space="preserve"> tasks. That</w:t></w:r>' + '<w:r><w:t>,</w:t></w:r>' + '<w:r><w:t xml:space="preserve"> plus the </w:t></w:r>' + '<w:proofErr w:type="spellStart"/>' + '<w:r><w:t>github</w:t></w:r>' + '<w:proofErr w:type="spellEnd"/>' + '<w:r><w:t xml:space="preserve"> creating it</w:t></w:r>' + '<w:r><w:t>s own complications, such as </w:t></w:r>' + '<w:proofErr w:type="spellStart"/>' + '<w:r><w:t>Khe</w:t></w:r>' + '<w:r><w:t>lil</w:t></w:r>' + '<w:proofErr w:type="spellEnd"/>' + '<w:r><w:t xml:space="preserve"> not being able to push at times when we tried to conjoin both of our parts together, which made debugging even harder</w:t></w:r>' + '<w:r w:rsidR="001A3AD1"><w:t>.</w:t></w:r>' + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Replace-ParagraphXml 9 $para9

Write-Host "All edits applied."
